# Commit: "done ! this works theres one type in 04/15/2009 that needs to be
# fixed mannually." — strips the footnote markers like " [1]", " [2]", ...
# " [6]" from the vaccine-name / brand-name cells across every worksheet, and
# collapses the embedded line-breaks in those same cells (e.g. "Hepatitis B
# [5]\nPediatric/Adolescent", "Recombivax\nHB") down to a single line by
# turning the newline into a plain space.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $ur = $ws.UsedRange
    $firstRow = $ur.Row
    $firstCol = $ur.Column
    $lastRow = $firstRow + $ur.Rows.Count - 1
    $lastCol = $firstCol + $ur.Columns.Count - 1

    for ($r = $firstRow; $r -le $lastRow; $r++) {
        for ($c = $firstCol; $c -le $lastCol; $c++) {
            $cell = $ws.Cells.Item($r, $c)
            $v = $cell.Value2

            if ($v -eq $null) { continue }
            if ($v.GetType().Name -ne "String") { continue }
            if ($v -notmatch "\[\d+\]" -and $v -notmatch "`n") { continue }

            $newVal = $v -replace "\[\d+\]", ""
            $newVal = $newVal.Replace("`n", " ")

            if ($newVal -ne $v) {
                $cell.Value = $newVal
            }
        }
    }
}
